$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    3  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    4  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 4.429675500412797 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    6  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 3.781711156805759 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    8  = @{ B = 0.3048080303191223; C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732; G = 5.014808316549482 }
    9  = @{ B = 0.3048080303191223; C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732; G = 5.014808316549482 }
    10 = @{ B = 0.01514828764759746; C = 0.002777888934908601; D = 3.900430680208489; E = 0.496779210170732; G = 4.415136066961727 }
    11 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    12 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    13 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 7.524616544037286 }
    14 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    15 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    16 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    17 = @{ B = 0.6753301551942219; C = 0.04240448674262143; D = 3.900430680208489;  E = 0.496779210170732; G = 5.114944532316064 }
    18 = @{ B = 0.3048080303191223; C = 0.3127903958511391;  D = 0.1575252929769615; E = 0.496779210170732; G = 1.271902929317955 }
    19 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    20 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    21 = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 0.8054896365839992; E = 0.496779210170732; G = 2.290389397800092 }
    22 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 2.997429241610044 }
    23 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    24 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
}
